$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column I ("Plan") before the existing formula column (which
# shifts from J to K).
$ws.Columns.Item(9).Insert()

# Header
$ws.Range("I1").Value = "Plan"

# Plan flag: 1 for Beijing (row 2), 0 for every other city (rows 3-37).
$ws.Range("I2").Value = 1
for ($r = 3; $r -le 37; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
}

# Rebuild the JSON-row formula in column K so it also emits the new
# "plan" field, sourced from the new column I.
$formulaRow2 = '="{ " & LOWER($A$1) & ": """ & A2 & """, " & LOWER($B$1) & ": " & B2 & ", " & LOWER($C$1) & ": " & C2 & ", " & LOWER($D$1) & ": """ & D2 & """, " & LOWER($E$1) & ": """ & E2 & """, " & LOWER($F$1) & ": """ & F2 & """, " & LOWER($G$1) & ": """ & G2 & """, " & LOWER($H$1) & ": """ & H2 & """, " & LOWER($I$1) & ": """ & I2 & """" &  " },"'
$ws.Range("K2").Formula = $formulaRow2

$formulaRow3 = '="{ " & LOWER($A$1) & ": """ & A3 & """, " & LOWER($B$1) & ": " & B3 & ", " & LOWER($C$1) & ": " & C3 & ", " & LOWER($D$1) & ": """ & D3 & """, " & LOWER($E$1) & ": """ & E3 & """, " & LOWER($F$1) & ": """ & F3 & """, " & LOWER($G$1) & ": """ & G3 & """, " & LOWER($H$1) & ": """ & H3 & """, " & LOWER($I$1) & ": """ & I3 & """" &  " },"'
$ws.Range("K3:K37").Formula = $formulaRow3

# Match the saved selection state (K2:K37).
$ws.Range("K2:K37").Select()
